# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" sheet (fund-holdings detail, same shape as the
# existing "2021-Q3"/"2021-Q4" sheets) right before the "总计" (totals)
# sheet, and prepends a matching summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new sheet and place it immediately before "总计".
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$newTmp = $wb.Worksheets.Add()
$newTmp.Name = "2022-Q1"

$totals = $wb.Worksheets.Item("总计")
$newTmp.Move([ref]$totals)

# Re-resolve by name - references captured before the Move/rename can be
# left pointing at the wrong sheet once indices shift.
$q1 = $wb.Worksheets.Item("2022-Q1")
$totals = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 2. Populate "2022-Q1" with the fund-holdings table.
#    Columns D-G hold numeric-looking figures that must stay TEXT
#    (e.g. fund codes like "011458" would lose their leading zero if
#    Excel auto-coerced them to numbers), so we pre-format B:G as Text.
# ---------------------------------------------------------------------
$q1.Range("B1:G4").NumberFormat = "@"

# Header row - copy the bold/centered/bordered header style from 总计.
$totals.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Index column (A) - copy the bold/centered/bordered style used by 总计's
# index column so the new sheet matches the existing ones.
$totals.Range("A2").Copy()
$q1.Range("A2:A4").PasteSpecial(-4122)

$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "400032"
$q1.Range("C2").Value = "东方主题精选混合"
$q1.Range("D2").Value = "29.31"
$q1.Range("E2").Value = "85.31"
$q1.Range("F2").Value = "2.99"
$q1.Range("G2").Value = "0.8764"
$q1.Range("H2").Value = 8

$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "011458"
$q1.Range("C3").Value = "东方鑫享价值成长一年持有期混合型证券投资基金A"
$q1.Range("D3").Value = "4.69"
$q1.Range("E3").Value = "80.98"
$q1.Range("F3").Value = "2.72"
$q1.Range("G3").Value = "0.1276"
$q1.Range("H3").Value = 10

$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "011459"
$q1.Range("C4").Value = "东方鑫享价值成长一年持有期混合型证券投资基金C"
$q1.Range("D4").Value = "1.71"
$q1.Range("E4").Value = "80.98"
$q1.Range("F4").Value = "2.72"
$q1.Range("G4").Value = "0.0465"
$q1.Range("H4").Value = 10

# H column (仓位排名) holds real numbers in the sibling sheets - restore
# General format there since the blanket Text format above only targeted
# B:G.

# ---------------------------------------------------------------------
# 3. Prepend the "2022-Q1" summary row to "总计", shifting the existing
#    2021-Q4 / 2021-Q3 rows down by one and renumbering the index column.
# ---------------------------------------------------------------------
$oldRow2Date = $totals.Range("B2").Value()
$oldRow2Count = $totals.Range("C2").Value()
$oldRow2Value = $totals.Range("D2").Value()

$oldRow3Date = $totals.Range("B3").Value()
$oldRow3Count = $totals.Range("C3").Value()
$oldRow3Value = $totals.Range("D3").Value()

# Extend the styled index column down to the new row 4.
$totals.Range("A2").Copy()
$totals.Range("A4").PasteSpecial(-4122)

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 3
$totals.Range("D2").Value = 1.05

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = $oldRow2Date
$totals.Range("C3").Value = $oldRow2Count
$totals.Range("D3").Value = $oldRow2Value

$totals.Range("A4").Value = 2
$totals.Range("B4").Value = $oldRow3Date
$totals.Range("C4").Value = $oldRow3Count
$totals.Range("D4").Value = $oldRow3Value
